$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Jd_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"

# Apply the same header formatting (bold font, thin border, centered/top) used by A1 to D1:E1
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# New job posting row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "ML Engineer/Data Scientist"
$ws.Range("C2").Value = "Please find the Job Description (JD) below for your reference"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2

# A2 previously held an (empty) wrap-text styled cell; clear that formatting back to default
$ws.Range("A2").ClearFormats()

# Reset the view selection/scroll position left over from the previous edit session
[void]$ws.Range("A1").Select()
